# Update Inscricoes summary counts (ResumoInscricoes) to reflect one
# additional registration ("Inscritos") on a handful of rows, cascading
# into the dependent Pagos / Inscrições homologadas totals where the
# diff shows them moving too.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# row 18: Inscritos 102 -> 103
$ws.Range("E18").Value = 103

# row 32: Inscritos 19 -> 20
$ws.Range("E32").Value = 20

# row 35: Inscritos 5 -> 6
$ws.Range("E35").Value = 6

# row 37: Inscritos 44 -> 45
$ws.Range("E37").Value = 45

# row 50: Inscritos 22 -> 23, Pagos 5 -> 6, Inscricoes homologadas 5 -> 6
$ws.Range("E50").Value = 23
$ws.Range("F50").Value = 6
$ws.Range("H50").Value = 6

# row 51: Inscritos 7 -> 8
$ws.Range("E51").Value = 8

# row 63: Inscritos 29 -> 30, Pagos 9 -> 10, Inscricoes homologadas 9 -> 10
$ws.Range("E63").Value = 30
$ws.Range("F63").Value = 10
$ws.Range("H63").Value = 10

# row 87: Inscritos 12 -> 13
$ws.Range("E87").Value = 13
